$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "87.557.52"
$ws.Range("E2").Value = "  -1.14%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.166.61"
$ws.Range("E3").Value = "  -6.81%  "

$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.86"
$ws.Range("E5").Value = "  -5.85%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "610.43"
$ws.Range("E6").Value = "  -5.94%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.380"
$ws.Range("E7").Value = "  -10.14%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.664"
$ws.Range("E8").Value = "  -1.67%  "

$ws.Range("E9").Value = "  -0.07%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.163.26"
$ws.Range("E10").Value = "  -6.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.534"
$ws.Range("E11").Value = "  -14.81%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.176"
$ws.Range("E12").Value = "  +3.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000240"
$ws.Range("E13").Value = "  -16.98%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.742.18"
$ws.Range("E14").Value = "  -7.39%  "

$ws.Range("E15").Value = "  -6.38%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.228.86"
$ws.Range("E16").Value = "  -1.48%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "32.18"
$ws.Range("E17").Value = "  -12.81%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.152.91"
$ws.Range("E18").Value = "  -7.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.06"
$ws.Range("E19").Value = "  -0.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.39"
$ws.Range("E20").Value = "  -11.13%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "414.68"
$ws.Range("E21").Value = "  -9.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.43"
$ws.Range("E22").Value = "  -12.88%  "

$ws.Range("E23").Value = "  -11.34%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.17"
$ws.Range("E24").Value = "  -6.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.91"
$ws.Range("E25").Value = "  -8.54%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.329.05"
$ws.Range("E26").Value = "  -7.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "73.39"
$ws.Range("E27").Value = "  -9.52%  "

$ws.Range("E28").Value = "  -10.86%  "

$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.26%  "

$ws.Range("B31").Value = "Cronos"
$ws.Range("C31").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.158"
$ws.Range("E31").Value = "  -17.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "542.31"
$ws.Range("E32").Value = "  -7.12%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.19"
$ws.Range("E33").Value = "  -13.59%  "

$ws.Range("E34").Value = "  -17.04%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.71"
$ws.Range("E35").Value = "  -9.52%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.85"
$ws.Range("E36").Value = "  -12.97%  "

$ws.Range("E37").Value = "  -8.74%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "21.74"
$ws.Range("E38").Value = "  -9.17%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "21.80"
$ws.Range("E39").Value = "  -0.22%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  +0.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.98"
$ws.Range("E41").Value = "  -5.82%  "

$ws.Range("E42").Value = "  +0.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.90"
$ws.Range("E43").Value = "  -10.47%  "

$ws.Range("E44").Value = "  -16.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "148.62"
$ws.Range("E45").Value = "  -6.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "172.38"
$ws.Range("E46").Value = "  -9.43%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "42.99"
$ws.Range("E47").Value = "  -8.11%  "

$ws.Range("E48").Value = "  -1.57%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.22"
$ws.Range("E49").Value = "  -14.61%  "

$ws.Range("E50").Value = "  -13.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.697"
$ws.Range("E51").Value = "  -12.08%  "
